$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7 (CIFF_1..CIFF_6): populate source_schema_path (column G) with the
# new JSON schema file path instead of "NOT APPL"
$schemaPath = "E:\ETL Automation\ETL_Framework\schema_files\Contact_info_schema.json"
foreach ($r in 2..7) {
    $ws.Range("G$r").Value = $schemaPath
}

# Row 8 (CIFF_7 / data_compare): switch source/target from the csv files to
# the new singleline.json file, and update source_type/target_type to "json"
$ws.Range("D8").Value = "json"
$ws.Range("C8").Value = "E:\ETL Automation\ETL_Framework\files\singleline.json"
$ws.Range("I8").Value = "json"
$ws.Range("H8").Value = "E:\ETL Automation\ETL_Framework\files\singleline.json"

# Update the view: active cell / selection on the active sheet
$ws.Range("C8").Select()
